# "changed siteurl for failing" - add a new UAT user row (test@click2cloud.com)
# to the Users sheet, mirroring the existing rows' layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New row 4: email (plain text, no hyperlink), password "123" (quote-prefixed
# text, like B2/B3), role "Test Manager".
$ws.Range("A4").Value = "test@click2cloud.com"

# Write B4 as literal text "123" first (leading apostrophe forces text so it
# doesn't get stored as a number), then copy B2's formatting (quote-prefix
# style) onto it so it matches B2/B3 exactly.
$ws.Range("B4").Value = "'123"
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Range("C4").Value = "Test Manager"

# Widen column C slightly and drop its "best fit" auto-sizing now that it
# holds a fixed custom width.
$ws.Columns.Item(3).ColumnWidth = 12.75

# Move the active selection, as recorded in the saved view state.
$ws.Range("B13").Select()
